$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find a paragraph by the exact text it contains and return the
# Paragraph object (via Word's Find on a throw-away Range so the document's
# selection/cursor state is left untouched).
# ---------------------------------------------------------------------------
function Find-ParagraphByText($text) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Find-ParagraphByText: text not found: $text"
    }
    return $rng.Paragraphs(1)
}

# ---------------------------------------------------------------------------
# Helper: toggle the presence of <w:lastRenderedPageBreak/> as the first
# child of the (sole) run that carries the given paragraph's text, by
# rewriting the paragraph's run content via Range.InsertXML (which replaces
# - rather than inserts into - the addressed range, and keeps the host
# <w:p>'s own attributes/paraId untouched).
# ---------------------------------------------------------------------------
function Set-LastRenderedPageBreak($paragraphText, [bool]$present) {
    $para = Find-ParagraphByText($paragraphText)
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $raw = $r.Text

    $escaped = $raw.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    if ($present) {
        $runInner = "<w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">$escaped</w:t>"
    } else {
        $runInner = "<w:t xml:space=`"preserve`">$escaped</w:t>"
    }

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1. Introduction paragraph: the sentence was split across two runs purely
#    because of an editing artifact ("...reduces human " / "error, and ...").
#    Collapse it back down to a single run with the identical text.
# ---------------------------------------------------------------------------
$introPara = Find-ParagraphByText("reduces human error, and improves the overall user experience.")
$introRange = $d.Range($introPara.Range.Start, $introPara.Range.End - 1)
$introText = $introRange.Text
$introRange.Delete()
$introRange.InsertAfter($introText)

# ---------------------------------------------------------------------------
# 2. Remove the "Frontend Technologies" sub-section entirely: the heading
#    line, its two bullet lines, and the blank spacer paragraph that
#    followed them.
# ---------------------------------------------------------------------------
$startRng = $d.Content
$startRng.Find.Execute("• Frontend Technologies:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("  - JavaScript for adding interactivity and dynamic content updates.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPara = $endRng.Paragraphs(1)
$blankPara = $endPara.Next()
$endPos = $blankPara.Range.End

$d.Range($startPos, $endPos).Delete()

# ---------------------------------------------------------------------------
# 3. Shift the rendered-page-break markers to match the new pagination that
#    results from the removed content above:
#      - drop it from the MySQL bullet (Database sub-section)
#      - add it to the "5. System Modules" heading
#      - drop it from the "6. Roles and Responsibilities" heading
#      - drop it from the "User Acceptance Testing" bullet
#      - add it to the "Handling concurrent bookings" bullet
# ---------------------------------------------------------------------------
Set-LastRenderedPageBreak("  - MySQL for storing user data, event details, booking records, and admin logs.") $false
Set-LastRenderedPageBreak("5. System Modules") $true
Set-LastRenderedPageBreak("6. Roles and Responsibilities") $false
Set-LastRenderedPageBreak("• User Acceptance Testing: Feedback from sample users was incorporated to enhance usability.") $false
Set-LastRenderedPageBreak("• Handling concurrent bookings and avoiding double bookings.") $true

Write-Output "Edit complete"
